$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2 through 51 (column A: Actual Consumption (MW), column B: Timestamp serial)
$data = @(
    @(6035, 45736),
    @(5957, 45736.01041666666),
    @(5972, 45736.02083333334),
    @(5979, 45736.03125),
    @(5907, 45736.04166666666),
    @(5857, 45736.05208333334),
    @(5863, 45736.0625),
    @(5884, 45736.07291666666),
    @(5850, 45736.08333333334),
    @(5876, 45736.09375),
    @(5854, 45736.10416666666),
    @(5868, 45736.11458333334),
    @(5829, 45736.125),
    @(5894, 45736.13541666666),
    @(5936, 45736.14583333334),
    @(6004, 45736.15625),
    @(6044, 45736.16666666666),
    @(6074, 45736.17708333334),
    @(6184, 45736.1875),
    @(6269, 45736.19791666666),
    @(6446, 45736.20833333334),
    @(6540, 45736.21875),
    @(6694, 45736.22916666666),
    @(6769, 45736.23958333334),
    @(6959, 45736.25),
    @(7017, 45736.26041666666),
    @(7027, 45736.27083333334),
    @(6990, 45736.28125),
    @(6926, 45736.29166666666),
    @(6869, 45736.30208333334),
    @(6791, 45736.3125),
    @(6592, 45736.32291666666),
    @(6459, 45736.33333333334),
    @(6344, 45736.34375),
    @(6241, 45736.35416666666),
    @(6084, 45736.36458333334),
    @(5945, 45736.375),
    @(5863, 45736.38541666666),
    @(5709, 45736.39583333334),
    @(5668, 45736.40625),
    @(5522, 45736.41666666666),
    @(5457, 45736.42708333334),
    @(5455, 45736.4375),
    @(5343, 45736.44791666666),
    @(5308, 45736.45833333334),
    @(5226, 45736.46875),
    @(5143, 45736.47916666666),
    @(5152, 45736.48958333334),
    @(5218, 45736.5),
    @(5258, 45736.51041666666)
)

$dateFormat = "YYYY-MM-DD HH:MM:SS"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $aValue = $data[$i][0]
    $bValue = $data[$i][1]

    $ws.Cells.Item($row, 1).Value = $aValue

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = $bValue
    $bCell.NumberFormat = $dateFormat
}
